# Slide 59, "Content Placeholder 2": the futures-pricing formula
#   FT=(S0+cost-benefit)*(1+rf)T
# needs "+cost-benefit)*(1+rf)" rewritten so that "cost" and "benefit"
# each get a subscript "0" (i.e. cost0 / benefit0), splitting the run
# into: "+cost" | "0"(sub) | "-benefit" | "0"(sub) | ")*(1+rf)"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(59)
$shape = $s.Shapes.Item("Content Placeholder 2")
$tr = $shape.TextFrame.TextRange

# Sanity check / starting text: "FT=(S0+cost-benefit)*(1+rf)T"

# Step 1: insert a "0" right after "cost" (before the "-" of "-benefit")
# by expanding the "-" at position 12 into "0-".
$dash = $tr.Characters(12, 1)
$dash.Text = "0-"

# Step 2: insert a "0" right after "benefit" (before the ")" that follows)
# by expanding the ")" now at position 21 into "0)".
$paren = $tr.Characters(21, 1)
$paren.Text = "0)"

# Step 3: re-assign the surrounding plain-text spans verbatim so they stay
# (or become) single runs rather than staying fragmented from the inserts.
$tr.Characters(13, 8).Text = "-benefit"
$tr.Characters(22, 8).Text = ")*(1+rf)"

# Step 4: make the two newly-inserted "0" characters subscript, matching
# the existing "0" / "T" baseline formatting used elsewhere in the formula.
$tr.Characters(12, 1).Font.BaselineOffset = -0.25
$tr.Characters(21, 1).Font.BaselineOffset = -0.25
